# Apply "average with safety stocks" edit.
#
# 1) On the "Productdata" sheet, the InventoryCosts (D), BackorderCosts (F)
#    and LostSale (I) columns for rows 2-11 are rescaled (divided by 2500).
# 2) On the "ForcastedStandardDeviation" sheet, the standard deviation
#    values (columns B-E) for rows 9-11 are zeroed out.

$wb = $excel.ActiveWorkbook

# --- Productdata sheet -----------------------------------------------
$ws = $wb.Worksheets.Item("Productdata")

$productData = @{
    2  = @{ D = 0.0016; F = 0.016; I = 0.16 }
    3  = @{ D = 0.0028; F = 0.028; I = 0.28 }
    4  = @{ D = 0.0024; F = 0.024; I = 0.24 }
    5  = @{ D = 0.0012; F = 0.012; I = 0.12 }
    6  = @{ D = 0.0012; F = 0.012; I = 0.12 }
    7  = @{ D = 0.0012; F = 0.012; I = 0.12 }
    8  = @{ D = 0.0008; F = 0.008; I = 0.08 }
    9  = @{ D = 0.0004; F = 0.004; I = 0.04 }
    10 = @{ D = 0.0004; F = 0.004; I = 0.04 }
    11 = @{ D = 0.0004; F = 0.004; I = 0.04 }
}

foreach ($row in $productData.Keys) {
    $vals = $productData[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("I$row").Value = $vals.I
}

# --- ForcastedStandardDeviation sheet ---------------------------------
$ws2 = $wb.Worksheets.Item("ForcastedStandardDeviation")

foreach ($row in 9..11) {
    $ws2.Range("B$row").Value = 0
    $ws2.Range("C$row").Value = 0
    $ws2.Range("D$row").Value = 0
    $ws2.Range("E$row").Value = 0
}
